$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / style) of the last existing data row (229)
# down onto the four new rows (230-233) so the new date cells keep the same
# date style (s="2") as the rest of column A, while B/C/D stay unstyled.
$srcRow = $ws.Range("A229:D229")
$srcRow.Copy()
$dstRow = $ws.Range("A230:D233")
$dstRow.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$data = @(
    @(230, 44304, 6, 95, 131.2716770993105),
    @(231, 44305, 10, 83, 114.6899915709765),
    @(232, 44306, 14, 79, 109.1627630615319),
    @(233, 44307, 2, 78, 107.7809559341707)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
